# Generate Report for handback
# Row for the "2a4ef3df-2a97-4c83-90a1-9918c8ba7ad3" source file transitions
# from "Ready for handoff" to "Handed back: in sync with en-US" on every
# sheet, and the "Latest Handback DateTime" column is refreshed for both
# rows on the per-language sheets.

$wb = $excel.ActiveWorkbook

$handedBack = "Handed back: in sync with en-US"

# --- Overview sheet --------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("B3").Value = $handedBack
$overview.Range("C3").Value = $handedBack

# --- zh-cn sheet -------------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("B3").Value = $handedBack
$zhcn.Range("G2").Value = "2016-01-26 06:31:32"
$zhcn.Range("G3").Value = "2016-01-26 06:31:32"

# --- de-de sheet -------------------------------------------------------
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("B3").Value = $handedBack
$dede.Range("G2").Value = "2016-01-26 06:31:56"
$dede.Range("G3").Value = "2016-01-26 06:31:56"
